$d = $word.ActiveDocument

# 1. Title "EID403: MACHINE LEARNING" - bump size from 14pt/28 half-points to 16pt/32 half-points
#    (covers both the title run and the trailing space run in paragraph 1, as well as the
#    paragraph mark's run properties)
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Font.Size = 16
$titlePara.Range.Font.SizeBi = 16

# 2. "L" + "  " + "T" -> single run "L  T"
$d.Content.Find.Execute("L  T", $true, $false, $false, $false, $false, $true, 1, $false, "L  T", 2) | Out-Null

# 3. "  " + "P" + "  " + "C " -> single run "  P  C "
$d.Content.Find.Execute("  P  C ", $true, $false, $false, $false, $false, $true, 1, $false, "  P  C ", 2) | Out-Null

# 4. "4" + "  " + "0" -> single run "4  0"
$d.Content.Find.Execute("4  0", $true, $false, $false, $false, $false, $true, 1, $false, "4  0", 2) | Out-Null

# 5. "  " + "0" + "  " + "4" -> single run "  0  4"
$d.Content.Find.Execute("  0  4", $true, $false, $false, $false, $false, $true, 1, $false, "  0  4", 2) | Out-Null

# 6. Module I body: merge all the separate sentence runs into one run
$d.Content.Find.Execute(" Learning problems, perspectives and issues, concept learning, version spaces and candidate eliminations, inductive bias, decision tree learning, representation, algorithm, heuristic space search", $true, $false, $false, $false, $false, $true, 1, $false, " Learning problems, perspectives and issues, concept learning, version spaces and candidate eliminations, inductive bias, decision tree learning, representation, algorithm, heuristic space search", 2) | Out-Null

# 7. "Module II " -> "Module II" (drop trailing plain-formatted space run)
$d.Content.Find.Execute("Module II ", $true, $false, $false, $false, $false, $true, 1, $false, "Module II", 2) | Out-Null

# 8. Neural Networks module body: merge all the separate sentence runs into one run
$d.Content.Find.Execute(" Neural network representation, problems, perceptron, multilayer networks and back propagation algorithms, advanced topics, genetic algorithms, hypothesis space search, genetic programming, models of evaluation and learning.", $true, $false, $false, $false, $false, $true, 1, $false, " Neural network representation, problems, perceptron, multilayer networks and back propagation algorithms, advanced topics, genetic algorithms, hypothesis space search, genetic programming, models of evaluation and learning.", 2) | Out-Null

# 9. "Bayes" + "theorem," -> single run " Bayes theorem,"
$d.Content.Find.Execute(" Bayes theorem,", $true, $false, $false, $false, $false, $true, 1, $false, " Bayes theorem,", 2) | Out-Null

# 10. "case-based" + " learning. " -> single run "case-based learning. "
$d.Content.Find.Execute("case-based learning. ", $true, $false, $false, $false, $false, $true, 1, $false, "case-based learning. ", 2) | Out-Null

# 11. Remove the _GoBack bookmark
$d.Bookmarks.Item("_GoBack").Delete()
